$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.402.92"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.607.30"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.69"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.11"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.068.58"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.329.47"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.616.18"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.94"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.50"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("E30").Value = "  +6.64%  "
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.31"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.05"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.835"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.812"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "274.71"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0956"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.948.86"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.29"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.41"
$ws.Range("E51").Value = "  -2.29%  "
